$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that vary row-to-row for this data block.
$cols = @("D", "J", "K", "L", "M", "P")

# Read the existing 20 rows (420-439) of data into memory first, since the
# shift operation needs the original values even after some of the
# destination rows are overwritten.
$source = @{}
for ($r = 420; $r -le 439; $r++) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Range($c + $r).Value2
    }
    $source[$r] = $rowVals
}

# The whole existing block (420-439) shifts down by two rows, to 422-441,
# making room for a brand-new pair of records at the top (420-421) and
# duplicating the final pair of records at the bottom (440-441).
for ($r = 439; $r -ge 420; $r--) {
    $dest = $r + 2
    foreach ($c in $cols) {
        $ws.Range($c + $dest).Value = $source[$r][$c]
    }
}

# New top record (Primera/Segunda pair) with fresh data.
$ws.Range("D420").Value = 44753
$ws.Range("J420").Value = 2400
$ws.Range("K420").Value = 600
$ws.Range("L420").Value = 700
$ws.Range("M420").Value = 650
$ws.Range("P420").Value = 325

$ws.Range("D421").Value = 44753
$ws.Range("J421").Value = 1400
$ws.Range("K421").Value = 500
$ws.Range("L421").Value = 550
$ws.Range("M421").Value = 525
$ws.Range("P421").Value = 262
